# Add a new "2022-Q1" sheet (positioned between "2021-Q2" and "总计")
# and record its totals as a new top row in the "总计" sheet.

function Set-TextCell {
    # Writes $value into $ws.Range($addr) forcing it to be stored as TEXT
    # (matters for values that look numeric, e.g. fund codes like "001417"
    # or decimal strings like "38.45"), then resets the cell's style back
    # to the sheet's default (no explicit style) by pasting formats from a
    # pristine, untouched scratch cell.
    param($ws, $addr, $value, $scratchAddr)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $ws.Range($scratchAddr).Copy()
    $r.PasteSpecial(-4122)
}

$wb = $excel.ActiveWorkbook

$sheetQ2 = $wb.Worksheets.Item("2021-Q2")

# --- Create the new "2022-Q1" sheet right after "2021-Q2" -------------
$sheetQ1 = $wb.Worksheets.Add($null, $sheetQ2)
$sheetQ1.Name = "2022-Q1"

# Match the page margins used by the other sheets in this workbook
# (0.75in / 1in / 0.5in, i.e. 54 / 72 / 36 points) instead of Excel's
# out-of-the-box defaults for a brand-new sheet.
$sheetQ1.PageSetup.LeftMargin = 54
$sheetQ1.PageSetup.RightMargin = 54
$sheetQ1.PageSetup.TopMargin = 72
$sheetQ1.PageSetup.BottomMargin = 72
$sheetQ1.PageSetup.HeaderMargin = 36
$sheetQ1.PageSetup.FooterMargin = 36

# Fetch the "总计" sheet *after* inserting the new sheet so the reference
# points at the correct (now shifted) worksheet rather than the new one.
$sheetTotal = $wb.Worksheets.Item("总计")

# Scratch cell (kept empty/untouched) used to reset styles back to default.
$scratch = "Z100"

# Copy the header-row style (bold/centered/bordered) from "总计" sheet
$sheetTotal.Range("B1:D1").Copy()
$sheetQ1.Range("B1:D1").PasteSpecial(-4122)
$sheetTotal.Range("B1").Copy()
$sheetQ1.Range("E1:H1").PasteSpecial(-4122)

# Copy the "index column" (A) style used in "总计" sheet for data rows
$sheetTotal.Range("A2").Copy()
$sheetQ1.Range("A2:A5").PasteSpecial(-4122)

# ---- Header row --------------------------------------------------------
$sheetQ1.Range("B1").Value = "基金代码"
$sheetQ1.Range("C1").Value = "基金名称"
$sheetQ1.Range("D1").Value = "基金规模"
$sheetQ1.Range("E1").Value = "股票总仓位"
$sheetQ1.Range("F1").Value = "仓位占比"
$sheetQ1.Range("G1").Value = "持有市值(亿元)"
$sheetQ1.Range("H1").Value = "仓位排名"

# ---- Data rows -----------------------------------------------------------
$sheetQ1.Range("A2").Value = 0
Set-TextCell $sheetQ1 "B2" "001417" $scratch
Set-TextCell $sheetQ1 "C2" "汇添富医疗服务灵活配置混合" $scratch
Set-TextCell $sheetQ1 "D2" "38.45" $scratch
Set-TextCell $sheetQ1 "E2" "77.97" $scratch
Set-TextCell $sheetQ1 "F2" "2.76" $scratch
Set-TextCell $sheetQ1 "G2" "1.0612" $scratch
$sheetQ1.Range("H2").Value = 10

$sheetQ1.Range("A3").Value = 1
Set-TextCell $sheetQ1 "B3" "015122" $scratch
Set-TextCell $sheetQ1 "C3" "汇添富医疗服务灵活配置混合D" $scratch
Set-TextCell $sheetQ1 "D3" "38.45" $scratch
Set-TextCell $sheetQ1 "E3" "77.97" $scratch
Set-TextCell $sheetQ1 "F3" "2.76" $scratch
Set-TextCell $sheetQ1 "G3" "1.0612" $scratch
$sheetQ1.Range("H3").Value = 10

$sheetQ1.Range("A4").Value = 2
Set-TextCell $sheetQ1 "B4" "010599" $scratch
Set-TextCell $sheetQ1 "C4" "汇添富高质量成长30一年持有期混合A" $scratch
Set-TextCell $sheetQ1 "D4" "20.34" $scratch
Set-TextCell $sheetQ1 "E4" "74.07" $scratch
Set-TextCell $sheetQ1 "F4" "2.53" $scratch
Set-TextCell $sheetQ1 "G4" "0.5146" $scratch
$sheetQ1.Range("H4").Value = 10

$sheetQ1.Range("A5").Value = 3
Set-TextCell $sheetQ1 "B5" "011259" $scratch
Set-TextCell $sheetQ1 "C5" "汇添富高质量成长30一年持有期混合C" $scratch
Set-TextCell $sheetQ1 "D5" "0.90" $scratch
Set-TextCell $sheetQ1 "E5" "74.07" $scratch
Set-TextCell $sheetQ1 "F5" "2.53" $scratch
Set-TextCell $sheetQ1 "G5" "0.0228" $scratch
$sheetQ1.Range("H5").Value = 10

$sheetQ1.Range($scratch).Clear()

# --- Insert a new summary row at the top of the "总计" sheet's data ----
$sheetTotal.Rows.Item(2).Insert()
$sheetTotal.Range("A2:D2").ClearFormats()

# Row 3 (the shifted original row) still carries the "A" column style; reuse it
$sheetTotal.Range("A3").Copy()
$sheetTotal.Range("A2").PasteSpecial(-4122)

$sheetTotal.Range("A2").Value = 0
$sheetTotal.Range("B2").Value = "2022-Q1"
$sheetTotal.Range("C2").Value = 4
$sheetTotal.Range("D2").Value = 2.66

# Renumber the shifted (originally first, index 0) row to its new index (1)
$sheetTotal.Range("A3").Value = 1
